$d = $word.ActiveDocument

# Locate the existing checklist bullet "Supprimer les méthodes en
# commentaires" — the new bullet about dropping redundant "this."
# prefixes belongs right after it (and before the blank paragraph /
# "Points importants :" heading that follows).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Supprimer les méthodes en commentaires*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 1) {
    $target = $d.Paragraphs.Item($targetIndex)

    # Insert a fresh paragraph right after it; Word clones the
    # paragraph formatting (style + bullet numbering) automatically.
    $target.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.Text = "Enlever les « this. » s’il n’y a pas d’ambiguïté au niveau des noms d’attributs"
}
